$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 11.57568926279069
$ws.Cells.Item(2, 3).Value = 8.411402818432375
$ws.Cells.Item(2, 4).Value = 5.779875983201483
$ws.Cells.Item(2, 5).Value = 12.33806647215485
$ws.Cells.Item(2, 6).Value = 27.57938297029872
$ws.Cells.Item(2, 8).Value = 7.344005520526261
$ws.Cells.Item(2, 9).Value = 24.88168758446652
$ws.Cells.Item(2, 11).Value = 8.990456234174884
$ws.Cells.Item(2, 12).Value = 10.13921418333481
$ws.Cells.Item(2, 13).Value = 13.85074040336108
$ws.Cells.Item(2, 14).Value = 20.01073722370112
$ws.Cells.Item(2, 15).Value = 24.77637013482199

$ws.Cells.Item(3, 2).Value = 11.35274404535034
$ws.Cells.Item(3, 3).Value = 8.368775653057321
$ws.Cells.Item(3, 4).Value = 5.726312211540698
$ws.Cells.Item(3, 5).Value = 12.36460280033805
$ws.Cells.Item(3, 6).Value = 27.61515176840143
$ws.Cells.Item(3, 8).Value = 7.344005520526261
$ws.Cells.Item(3, 9).Value = 24.96151525016302
$ws.Cells.Item(3, 11).Value = 8.824049441003369
$ws.Cells.Item(3, 12).Value = 10.14659296349595
$ws.Cells.Item(3, 13).Value = 13.81954460307029
$ws.Cells.Item(3, 14).Value = 20.0681507714159
$ws.Cells.Item(3, 15).Value = 24.84301352394822

$ws.Cells.Item(4, 2).Value = 11.21566733690699
$ws.Cells.Item(4, 3).Value = 8.342112319003272
$ws.Cells.Item(4, 4).Value = 5.692592784133669
$ws.Cells.Item(4, 5).Value = 12.38243338315666
$ws.Cells.Item(4, 6).Value = 27.6433379468749
$ws.Cells.Item(4, 8).Value = 7.344005520526261
$ws.Cells.Item(4, 9).Value = 25.01454087973661
$ws.Cells.Item(4, 11).Value = 8.721702006289293
$ws.Cells.Item(4, 12).Value = 10.15249578877338
$ws.Cells.Item(4, 13).Value = 13.80233469023047
$ws.Cells.Item(4, 14).Value = 20.10506738927893
$ws.Cells.Item(4, 15).Value = 24.88844642166259

$ws.Cells.Item(5, 2).Value = 11.15983759346538
$ws.Cells.Item(5, 3).Value = 8.331126388260564
$ws.Cells.Item(5, 4).Value = 5.678648164320137
$ws.Cells.Item(5, 5).Value = 12.39008638610724
$ws.Cells.Item(5, 6).Value = 27.65638779086833
$ws.Cells.Item(5, 8).Value = 7.344005520526261
$ws.Cells.Item(5, 9).Value = 25.03715771493107
$ws.Cells.Item(5, 11).Value = 8.680006464146386
$ws.Cells.Item(5, 12).Value = 10.15524689660552
$ws.Cells.Item(5, 13).Value = 13.79581564970956
$ws.Cells.Item(5, 14).Value = 20.12053095236785
$ws.Cells.Item(5, 15).Value = 24.90809458980645

$ws.Cells.Item(6, 2).Value = 11.15057121412949
$ws.Cells.Item(6, 3).Value = 8.329294993222106
$ws.Cells.Item(6, 4).Value = 5.67632052545126
$ws.Cells.Item(6, 5).Value = 12.39138053936731
$ws.Cells.Item(6, 6).Value = 27.65864910460508
$ws.Cells.Item(6, 8).Value = 7.344005520526261
$ws.Cells.Item(6, 9).Value = 25.04097412678662
$ws.Cells.Item(6, 11).Value = 8.673085305471083
$ws.Cells.Item(6, 12).Value = 10.15572460911251
$ws.Cells.Item(6, 13).Value = 13.79476315752914
$ws.Cells.Item(6, 14).Value = 20.12312405184801
$ws.Cells.Item(6, 15).Value = 24.91142560415426

$ws.Cells.Item(7, 2).Value = 11.21491417227256
$ws.Cells.Item(7, 3).Value = 8.341964642681456
$ws.Cells.Item(7, 4).Value = 5.692405539756382
$ws.Cells.Item(7, 5).Value = 12.3825350272269
$ws.Cells.Item(7, 6).Value = 27.64350761223974
$ws.Cells.Item(7, 8).Value = 7.344005520526261
$ws.Cells.Item(7, 9).Value = 25.01484181529224
$ws.Cells.Item(7, 11).Value = 8.72113956368014
$ws.Cells.Item(7, 12).Value = 10.15253149087058
$ws.Cells.Item(7, 13).Value = 13.80224476482389
$ws.Cells.Item(7, 14).Value = 20.10527423519199
$ws.Cells.Item(7, 15).Value = 24.88870681413931

$ws.Cells.Item(8, 2).Value = 11.49890249683496
$ws.Cells.Item(8, 3).Value = 8.396808381242572
$ws.Cells.Item(8, 4).Value = 5.76158282640305
$ws.Cells.Item(8, 5).Value = 12.34689737920923
$ws.Cells.Item(8, 6).Value = 27.59042362197627
$ws.Cells.Item(8, 8).Value = 7.344005520526261
$ws.Cells.Item(8, 9).Value = 24.90837955863756
$ws.Cells.Item(8, 11).Value = 8.933147954733686
$ws.Cells.Item(8, 12).Value = 10.14147396146905
$ws.Cells.Item(8, 13).Value = 13.83958374353295
$ws.Cells.Item(8, 14).Value = 20.0301886790376
$ws.Cells.Item(8, 15).Value = 24.79841135351477

$ws.Cells.Item(9, 2).Value = 12.05101497620483
$ws.Cells.Item(9, 3).Value = 8.500346187602444
$ws.Cells.Item(9, 4).Value = 5.890419064262623
$ws.Cells.Item(9, 5).Value = 12.28919261571202
$ws.Cells.Item(9, 6).Value = 27.53574918551126
$ws.Cells.Item(9, 8).Value = 7.344005520526261
$ws.Cells.Item(9, 9).Value = 24.73144106361068
$ws.Cells.Item(9, 11).Value = 9.345153287727014
$ws.Cells.Item(9, 12).Value = 10.13064981085868
$ws.Cells.Item(9, 13).Value = 13.92799132840898
$ws.Cells.Item(9, 14).Value = 19.89609731982577
$ws.Cells.Item(9, 15).Value = 24.65720150289088

$ws.Cells.Item(10, 2).Value = 12.44948469270484
$ws.Cells.Item(10, 3).Value = 8.573810382261948
$ws.Cells.Item(10, 4).Value = 5.980624122030435
$ws.Cells.Item(10, 5).Value = 12.25420025585147
$ws.Cells.Item(10, 6).Value = 27.5257312437812
$ws.Cells.Item(10, 8).Value = 7.344005520526261
$ws.Cells.Item(10, 9).Value = 24.6208604073553
$ws.Cells.Item(10, 11).Value = 9.642513765553016
$ws.Cells.Item(10, 12).Value = 10.12927665151722
$ws.Cells.Item(10, 13).Value = 14.00185749826657
$ws.Cells.Item(10, 14).Value = 19.80551850103056
$ws.Cells.Item(10, 15).Value = 24.57537437890565

$ws.Cells.Item(11, 2).Value = 12.62836285587469
$ws.Cells.Item(11, 3).Value = 8.606631694091629
$ws.Cells.Item(11, 4).Value = 6.020634394715628
$ws.Cells.Item(11, 5).Value = 12.2398839384368
$ws.Cells.Item(11, 6).Value = 27.52771327147611
$ws.Cells.Item(11, 8).Value = 7.344005520526261
$ws.Cells.Item(11, 9).Value = 24.57477242455027
$ws.Cells.Item(11, 11).Value = 9.776024249251392
$ws.Cells.Item(11, 12).Value = 10.13007038775871
$ws.Cells.Item(11, 13).Value = 14.03731621514002
$ws.Cells.Item(11, 14).Value = 19.76601863288465
$ws.Cells.Item(11, 15).Value = 24.54291948694539

$ws.Cells.Item(12, 2).Value = 12.69569073700127
$ws.Cells.Item(12, 3).Value = 8.618971370480773
$ws.Cells.Item(12, 4).Value = 6.035632995890084
$ws.Cells.Item(12, 5).Value = 12.23469265243305
$ws.Cells.Item(12, 6).Value = 27.52940236268616
$ws.Cells.Item(12, 8).Value = 7.344005520526261
$ws.Cells.Item(12, 9).Value = 24.5579265447453
$ws.Cells.Item(12, 11).Value = 9.826280686480796
$ws.Cells.Item(12, 12).Value = 10.13057390198584
$ws.Cells.Item(12, 13).Value = 14.05100315682215
$ws.Cells.Item(12, 14).Value = 19.75130497737662
$ws.Cells.Item(12, 15).Value = 24.53131610064587

$ws.Cells.Item(13, 2).Value = 12.68120969209687
$ws.Cells.Item(13, 3).Value = 8.616317809114515
$ws.Cells.Item(13, 4).Value = 6.032409637484744
$ws.Cells.Item(13, 5).Value = 12.23580046561025
$ws.Cells.Item(13, 6).Value = 27.52899688040958
$ws.Cells.Item(13, 8).Value = 7.344005520526261
$ws.Cells.Item(13, 9).Value = 24.56152761930299
$ws.Cells.Item(13, 11).Value = 9.815471189816416
$ws.Cells.Item(13, 12).Value = 10.13045645233575
$ws.Cells.Item(13, 13).Value = 14.04804401029011
$ws.Cells.Item(13, 14).Value = 19.75446298749262
$ws.Cells.Item(13, 15).Value = 24.53378454962517

$ws.Cells.Item(14, 2).Value = 12.63391053081436
$ws.Cells.Item(14, 3).Value = 8.60764868870492
$ws.Cells.Item(14, 4).Value = 6.021871422133715
$ws.Cells.Item(14, 5).Value = 12.23945224091159
$ws.Cells.Item(14, 6).Value = 27.52783343689536
$ws.Cells.Item(14, 8).Value = 7.344005520526261
$ws.Cells.Item(14, 9).Value = 24.57337434353972
$ws.Cells.Item(14, 11).Value = 9.780165176083647
$ws.Cells.Item(14, 12).Value = 10.13010775101992
$ws.Cells.Item(14, 13).Value = 14.03843709133942
$ws.Cells.Item(14, 14).Value = 19.7648032458323
$ws.Cells.Item(14, 15).Value = 24.54195110493188

$ws.Cells.Item(15, 2).Value = 12.60488321821591
$ws.Cells.Item(15, 3).Value = 8.60232691566272
$ws.Cells.Item(15, 4).Value = 6.015396462184384
$ws.Cells.Item(15, 5).Value = 12.24171899973128
$ws.Cells.Item(15, 6).Value = 27.52724295533509
$ws.Cells.Item(15, 8).Value = 7.344005520526261
$ws.Cells.Item(15, 9).Value = 24.58070981891753
$ws.Cells.Item(15, 11).Value = 9.758498634766934
$ws.Cells.Item(15, 12).Value = 10.12992055896488
$ws.Cells.Item(15, 13).Value = 14.03258614273712
$ws.Cells.Item(15, 14).Value = 19.77116870831139
$ws.Cells.Item(15, 15).Value = 24.54704279149088

$ws.Cells.Item(16, 2).Value = 12.43774080459659
$ws.Cells.Item(16, 3).Value = 8.571653138534348
$ws.Cells.Item(16, 4).Value = 5.977988302290963
$ws.Cells.Item(16, 5).Value = 12.25516806475383
$ws.Cells.Item(16, 6).Value = 27.52573316143254
$ws.Cells.Item(16, 8).Value = 7.344005520526261
$ws.Cells.Item(16, 9).Value = 24.62395723290168
$ws.Cells.Item(16, 11).Value = 9.633748960458508
$ws.Cells.Item(16, 12).Value = 10.12925323540448
$ws.Cells.Item(16, 13).Value = 13.99957691484025
$ws.Cells.Item(16, 14).Value = 19.80813412502891
$ws.Cells.Item(16, 15).Value = 24.5775914155775

$ws.Cells.Item(17, 2).Value = 12.33454460531684
$ws.Cells.Item(17, 3).Value = 8.552680444375687
$ws.Cells.Item(17, 4).Value = 5.954773542777228
$ws.Cells.Item(17, 5).Value = 12.26382866264432
$ws.Cells.Item(17, 6).Value = 27.52648073939828
$ws.Cells.Item(17, 8).Value = 7.344005520526261
$ws.Cells.Item(17, 9).Value = 24.65156817147299
$ws.Cells.Item(17, 11).Value = 9.556733273345991
$ws.Cells.Item(17, 12).Value = 10.12920643922667
$ws.Cells.Item(17, 13).Value = 13.97979712026147
$ws.Cells.Item(17, 14).Value = 19.8312471068739
$ws.Cells.Item(17, 15).Value = 24.59755402069115

$ws.Cells.Item(18, 2).Value = 12.27496785318845
$ws.Cells.Item(18, 3).Value = 8.541711755757433
$ws.Cells.Item(18, 4).Value = 5.941324896674723
$ws.Cells.Item(18, 5).Value = 12.26896080853218
$ws.Cells.Item(18, 6).Value = 27.52752637944092
$ws.Cells.Item(18, 8).Value = 7.344005520526261
$ws.Cells.Item(18, 9).Value = 24.66784603111692
$ws.Cells.Item(18, 11).Value = 9.512272906571333
$ws.Cells.Item(18, 12).Value = 10.12931309541837
$ws.Cells.Item(18, 13).Value = 13.96859552769376
$ws.Cells.Item(18, 14).Value = 19.84470161863103
$ws.Cells.Item(18, 15).Value = 24.60948479051302

$ws.Cells.Item(19, 2).Value = 12.25476030140931
$ws.Cells.Item(19, 3).Value = 8.537988405168649
$ws.Cells.Item(19, 4).Value = 5.936755063520695
$ws.Cells.Item(19, 5).Value = 12.27072437542132
$ws.Cells.Item(19, 6).Value = 27.5279861981237
$ws.Cells.Item(19, 8).Value = 7.344005520526261
$ws.Cells.Item(19, 9).Value = 24.67342556859754
$ws.Cells.Item(19, 11).Value = 9.497192902396275
$ws.Cells.Item(19, 12).Value = 10.12937217822951
$ws.Cells.Item(19, 13).Value = 13.96483317857603
$ws.Cells.Item(19, 14).Value = 19.84928469023228
$ws.Cells.Item(19, 15).Value = 24.61360139948962

$ws.Cells.Item(20, 2).Value = 12.34555338294733
$ws.Cells.Item(20, 3).Value = 8.554705948100933
$ws.Cells.Item(20, 4).Value = 5.957254780113281
$ws.Cells.Item(20, 5).Value = 12.26289112261379
$ws.Cells.Item(20, 6).Value = 27.52633745158517
$ws.Cells.Item(20, 8).Value = 7.344005520526261
$ws.Cells.Item(20, 9).Value = 24.64858787176763
$ws.Cells.Item(20, 11).Value = 9.564948948323988
$ws.Cells.Item(20, 12).Value = 10.12919760367449
$ws.Cells.Item(20, 13).Value = 13.98188463159186
$ws.Cells.Item(20, 14).Value = 19.82877008344163
$ws.Cells.Item(20, 15).Value = 24.59538251089276

$ws.Cells.Item(21, 2).Value = 12.6478150593449
$ws.Cells.Item(21, 3).Value = 8.610197460400556
$ws.Cells.Item(21, 4).Value = 6.024970928064588
$ws.Cells.Item(21, 5).Value = 12.23837338696554
$ws.Cells.Item(21, 6).Value = 27.52814971430247
$ws.Cells.Item(21, 8).Value = 7.344005520526261
$ws.Cells.Item(21, 9).Value = 24.56987820612718
$ws.Cells.Item(21, 11).Value = 9.790543948604363
$ws.Cells.Item(21, 12).Value = 10.13020467375548
$ws.Cells.Item(21, 13).Value = 14.04125189482565
$ws.Cells.Item(21, 14).Value = 19.76175944535382
$ws.Cells.Item(21, 15).Value = 24.53953375003594

$ws.Cells.Item(22, 2).Value = 12.84294634737516
$ws.Cells.Item(22, 3).Value = 8.645943959608452
$ws.Cells.Item(22, 4).Value = 6.068336654258331
$ws.Cells.Item(22, 5).Value = 12.2236900264654
$ws.Cells.Item(22, 6).Value = 27.53480336634743
$ws.Cells.Item(22, 8).Value = 7.344005520526261
$ws.Cells.Item(22, 9).Value = 24.52197303457419
$ws.Cells.Item(22, 11).Value = 9.936208007694272
$ws.Cells.Item(22, 12).Value = 10.13204523118389
$ws.Cells.Item(22, 13).Value = 14.08156069007261
$ws.Cells.Item(22, 14).Value = 19.71938635626625
$ws.Cells.Item(22, 15).Value = 24.50703555739896

$ws.Cells.Item(23, 2).Value = 12.73904283088686
$ws.Cells.Item(23, 3).Value = 8.626914004594752
$ws.Cells.Item(23, 4).Value = 6.045274697716795
$ws.Cells.Item(23, 5).Value = 12.23140428758428
$ws.Cells.Item(23, 6).Value = 27.53075251474521
$ws.Cells.Item(23, 8).Value = 7.344005520526261
$ws.Cells.Item(23, 9).Value = 24.54721723150647
$ws.Cells.Item(23, 11).Value = 9.858641955212432
$ws.Cells.Item(23, 12).Value = 10.13095506149635
$ws.Cells.Item(23, 13).Value = 14.05991161629614
$ws.Cells.Item(23, 14).Value = 19.74187188874492
$ws.Cells.Item(23, 15).Value = 24.52401400486011

$ws.Cells.Item(24, 2).Value = 12.34057708161124
$ws.Cells.Item(24, 3).Value = 8.553790407755452
$ws.Cells.Item(24, 4).Value = 5.956133330566979
$ws.Cells.Item(24, 5).Value = 12.26331450769846
$ws.Cells.Item(24, 6).Value = 27.52640031356381
$ws.Cells.Item(24, 8).Value = 7.344005520526261
$ws.Cells.Item(24, 9).Value = 24.64993400697641
$ws.Cells.Item(24, 11).Value = 9.561235208060253
$ws.Cells.Item(24, 12).Value = 10.12920118213006
$ws.Cells.Item(24, 13).Value = 13.98094033739099
$ws.Cells.Item(24, 14).Value = 19.82988942685225
$ws.Cells.Item(24, 15).Value = 24.59636283634959

$ws.Cells.Item(25, 2).Value = 11.90261629527803
$ws.Cells.Item(25, 3).Value = 8.472782361473815
$ws.Cells.Item(25, 4).Value = 5.856326282879941
$ws.Cells.Item(25, 5).Value = 12.30350151989022
$ws.Cells.Item(25, 6).Value = 27.54524434100038
$ws.Cells.Item(25, 8).Value = 7.344005520526261
$ws.Cells.Item(25, 9).Value = 24.77589860881516
$ws.Cells.Item(25, 11).Value = 9.234418912979219
$ws.Cells.Item(25, 12).Value = 10.13241935430583
$ws.Cells.Item(25, 13).Value = 13.90248317618695
$ws.Cells.Item(25, 14).Value = 24.6915570317898
